$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AutoTestAdmin"
$ws.Range("B2").Value = "AutoTestUser"
